$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently sits alone in
#    the paragraph right after "README". Its start/end collapse onto the
#    same spot, so removing it just leaves an empty paragraph behind
#    (matching the target "<w:p/>").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Re-insert the "_GoBack" bookmark, collapsed, right after the
#    "spring-boot:run" text (but before the trailing proofErr/paragraph
#    mark) inside the "mvn spring-boot:run" paragraph.
#
#    The engine's Bookmarks.Add mishandles a *degenerate* range that
#    lands exactly on a paragraph-mark position (it silently resets to
#    an unrelated range). Work around this by temporarily inserting a
#    one-character placeholder at the target spot, anchoring the
#    bookmark to that (non-degenerate) character range, and then
#    deleting the placeholder again -- the bookmark collapses correctly
#    onto the original spot once the placeholder text is gone.
#    InsertBefore (rather than InsertAfter) keeps the placeholder, and
#    therefore the final bookmark, ordered ahead of the zero-width
#    "spellEnd" proofErr marker already sitting at that position.
$find = $d.Content
[void]$find.Find.Execute("spring-boot:run", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find.Collapse(0)
$find.InsertBefore("X")

$placeholder = $d.Duplicate
$placeholder.Start = $find.Start
$placeholder.End = $find.Start + 1
$d.Bookmarks.Add("_GoBack", $placeholder)

$placeholder.Text = ""

# 3. Remove the now-superfluous empty paragraph that used to follow the
#    "mvn spring-boot:run" paragraph. Deleting the *whole* paragraph's own
#    Range (rather than reaching into it from the previous paragraph) keeps
#    the surrounding paragraphs' own pPr/formatting untouched -- merging
#    across a paragraph mark instead would make the previous paragraph
#    inherit the (blank) formatting of the deleted one.
$again = $d.Content
[void]$again.Find.Execute("spring-boot:run", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraIndex = $d.Range(0, $again.Start).Paragraphs.Count
$empty = $d.Paragraphs.Item($paraIndex + 1)
$empty.Range.Delete()
